# Commit: "optimized experience type = long -> int"
#
# The underlying content change is a rename of two header cells in the
# quests sheet:
#   D1: "npcKillAmountNeeded"  -> "npcAmountNeed"
#   E1: "userKillAmountNeeded" -> "userAmountNeed"
#
# Renaming these (previously mid-table) shared strings pushes them to the
# end of the shared-string table, which is why every other shared-string
# index in the sheet appears to shift in the diff -- the underlying text
# for every other cell (B2:C14, H1, etc.) is unchanged.
#
# Along with the rename, the header column autosizes to the new (shorter)
# text, the active selection moved to J4, and columns D & E (previously a
# single merged-width col entry) now carry their own individual widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two headers -- this alone reshuffles the shared-string table
# exactly like the diff (the two renamed strings move to the tail of the
# table, shifting every later index down by 2).
$ws.Range("D1").Value = "npcAmountNeed"
$ws.Range("E1").Value = "userAmountNeed"

# Reflect the narrower bestFit column widths that Excel recomputed after
# the header text got shorter.
$ws.Range("D1").EntireColumn.ColumnWidth = 17.0
$ws.Range("E1").EntireColumn.ColumnWidth = 17.65

# Move the active selection to match the post-edit cursor position.
$ws.Range("J4").Select()
